$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6272
$ws.Range("C3").Value = 0.5898
$ws.Range("C4").Value = 0.5433
$ws.Range("C5").Value = 0.3398
$ws.Range("C6").Value = 0.029
